$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 - this shifts the existing rows 19..131
# down to 20..132 (and the sheet's UsedRange/dimension grows to R132).
$ws.Rows.Item(19).Insert()

# Populate the freshly inserted row 19 with the new weekly data point.
# Columns A, B, C, E, F, G, N, Q, R are constant for every data row in this
# sheet (same market/category), so copy them down from the row below
# (which now holds what used to be row 19's data).
$ws.Cells.Item(19, 1).Value = $ws.Cells.Item(20, 1).Value()
$ws.Cells.Item(19, 2).Value = $ws.Cells.Item(20, 2).Value()
$ws.Cells.Item(19, 3).Value = $ws.Cells.Item(20, 3).Value()
$ws.Cells.Item(19, 4).Value = 44613
$ws.Cells.Item(19, 5).Value = $ws.Cells.Item(20, 5).Value()
$ws.Cells.Item(19, 6).Value = $ws.Cells.Item(20, 6).Value()
$ws.Cells.Item(19, 7).Value = $ws.Cells.Item(20, 7).Value()
$ws.Cells.Item(19, 8).Value = "Camote"
$ws.Cells.Item(19, 9).Value = "1a (cosecha)"
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 350
$ws.Cells.Item(19, 12).Value = 400
$ws.Cells.Item(19, 13).Value = 375
$ws.Cells.Item(19, 14).Value = $ws.Cells.Item(20, 14).Value()
$ws.Cells.Item(19, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 16).Value = 375
$ws.Cells.Item(19, 17).Value = $ws.Cells.Item(20, 17).Value()
$ws.Cells.Item(19, 18).Value = $ws.Cells.Item(20, 18).Value()

# Give the new date cell the same number format as the rest of column D.
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat()
